$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-124 down to 37-125
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with its data
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 45274
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112030
$ws.Range("G36").Value = "Poroto granado"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 80
$ws.Range("K36").Value = 60000
$ws.Range("L36").Value = 60000
$ws.Range("M36").Value = 60000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 2400
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
